$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Ativação: date changed 01/01/2020 -> 01/01/2022
# Force text format so Excel doesn't auto-convert the date string to a date serial number
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "01/01/2022"

# Row 14: Programa resumido
$ws.Range("B14").Value = "Barragens e Reservatórios. Usos da água demandados para o interesse humano e Panorama Geral da Engenharia dos Recursos Hídricos. Impactos Ambientais dos Usos da Água. Gestão dos Recursos Hídricos."
$ws.Range("C14").Value = "Barragens e Reservatórios. Usos da água demandados para o interesse humano e Panorama Geral da Engenharia dos Recursos Hídricos. Impactos Ambientais dos Usos da Água. Gestão dos Recursos Hídricos."

# Row 15: Short syllabus
$ws.Range("B15").Value = "Barrages and Water Tanks. The use of the water demanded for the human interest and Hydric Resources Engineering General View. The Usage of the Water Enviropnmetal Impacts. Hydric Resources Management."
$ws.Range("C15").Value = "Barrages and Water Tanks. The use of the water demanded for the human interest and Hydric Resources Engineering General View. The Usage of the Water Enviropnmetal Impacts. Hydric Resources Management."

# Row 16: Programa
$ws.Range("B16").Value = "- Políticas Públicas, Balanço Hídrico,- Demanda de água e disponibilidade dos recursos hídricos: Abastecimento Humano, águas para Agropecuária e indústria. - Hidreletricidade. - Barragens e Reservatórios,- Navegação Interior.- Águas Subterrâneas.- Gerenciamento dos Recursos Hídricos.- Hidroeconomia- Relação entre saneamento e qualidade da água"
$ws.Range("C16").Value = "- Políticas Públicas, Balanço Hídrico,- Demanda de água e disponibilidade dos recursos hídricos: Abastecimento Humano, águas para Agropecuária e indústria. - Hidreletricidade. - Barragens e Reservatórios,- Navegação Interior.- Águas Subterrâneas.- Gerenciamento dos Recursos Hídricos.- Hidroeconomia- Relação entre saneamento e qualidade da água"

# Row 17: Syllabus
$ws.Range("B17").Value = "- Public Policies, Water Balance,- Water demand and availability of water resources: Human Supply, water for Agriculture and industry.- Hydroelectricity.- Dams and Reservoirs,- Inland navigation.- Groundwater.- Water Resources Management.- Hydroeconomics- Relationship between sanitation and water quality"
$ws.Range("C17").Value = "- Public Policies, Water Balance,- Water demand and availability of water resources: Human Supply, water for Agriculture and industry.- Hydroelectricity.- Dams and Reservoirs,- Inland navigation.- Groundwater.- Water Resources Management.- Hydroeconomics- Relationship between sanitation and water quality"

# Row 19: Método
$ws.Range("B19").Value = "Avaliação baseada em trabalhos com dados reais, exercícios, trabalhos práticos e relatórios."
$ws.Range("C19").Value = "Avaliação baseada em trabalhos com dados reais, exercícios, trabalhos práticos e relatórios."

# Row 20: Critério
$ws.Range("B20").Value = "Média ponderada das notas atribuídas aos exercícios e trabalhos práticos e relatórios."
$ws.Range("C20").Value = "Média ponderada das notas atribuídas aos exercícios e trabalhos práticos e relatórios."
